$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.085.88"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.609.30"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'604.39"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").Value = "'145.41"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.585"
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").Value = "2.608.24"
$ws.Range("E9").Value = "  -1.48%  "
$ws.Range("E10").Value = "  +1.17%  "
$ws.Range("D11").Value = "'5.50"
$ws.Range("E11").Value = "  -2.92%  "
$ws.Range("E12").Value = "  +4.08%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "'27.12"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").Value = "3.076.16"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "62.920.60"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "'0.0000146"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").Value = "2.589.51"
$ws.Range("E18").Value = "  -2.34%  "
$ws.Range("D19").Value = "'11.43"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("E20").Value = "  +3.02%  "
$ws.Range("D21").Value = "'341.68"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "'6.84"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "'5.69"
$ws.Range("E24").Value = "  -1.29%  "
$ws.Range("D25").Value = "'66.09"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").Value = "'1.69"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").Value = "'1.59"
$ws.Range("E27").Value = "  +4.39%  "
$ws.Range("D28").Value = "'9.00"
$ws.Range("E28").Value = "  +6.47%  "
$ws.Range("D29").Value = "'559.56"
$ws.Range("E29").Value = "  +4.43%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "'0.161"
$ws.Range("E31").Value = "  -2.74%  "
$ws.Range("D32").Value = "'7.75"
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("D34").Value = "0.0₃0842"
$ws.Range("E34").Value = "  +4.19%  "
$ws.Range("E35").Value = "  -5.30%  "
$ws.Range("D36").Value = "'5.15"
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("D37").Value = "'167.81"
$ws.Range("E37").Value = "  -3.42%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").Value = "'1.92"
$ws.Range("E40").Value = "  +4.10%  "
$ws.Range("D41").Value = "'18.96"
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "'164.64"
$ws.Range("E43").Value = "  -4.65%  "
$ws.Range("D44").Value = "'39.51"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("D45").Value = "'3.74"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").Value = "'21.73"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D47").Value = "'0.0563"
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("D48").Value = "'0.623"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("D49").Value = "'0.0245"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").Value = "'0.0955"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").Value = "'1.91"
$ws.Range("E51").Value = "  +11.07%  "
